$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.419.59'
$ws.Range('E2').Value = '  -4.46%  '
$ws.Range('D3').Value = '3.309.69'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.06'
$ws.Range('E5').Value = '  -4.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.47'
$ws.Range('E6').Value = '  -5.92%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.599'
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('D9').Value = '3.309.91'
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.129'
$ws.Range('E10').Value = '  -3.84%  '
$ws.Range('E11').Value = '  -1.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.403'
$ws.Range('E12').Value = '  -4.62%  '
$ws.Range('D13').Value = '3.881.39'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.137'
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.99'
$ws.Range('E15').Value = '  -4.55%  '
$ws.Range('D16').Value = '66.509.63'
$ws.Range('E16').Value = '  -4.34%  '
$ws.Range('E17').Value = '  -3.43%  '
$ws.Range('D18').Value = '3.291.21'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '435.80'
$ws.Range('E19').Value = '  +3.27%  '
$ws.Range('E20').Value = '  -0.89%  '
$ws.Range('E21').Value = '  -2.96%  '
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.44'
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').Value = '3.448.01'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.512'
$ws.Range('E26').Value = '  -1.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000118'
$ws.Range('E27').Value = '  -2.89%  '
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.99'
$ws.Range('E29').Value = '  -7.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('E31').Value = '  -2.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '22.69'
$ws.Range('E32').Value = '  -1.88%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.29'
$ws.Range('E33').Value = '  -6.08%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.75'
$ws.Range('E35').Value = '  -4.23%  '
$ws.Range('E36').Value = '  -5.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.50'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '159.60'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.25'
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.84'
$ws.Range('E40').Value = '  -5.36%  '
$ws.Range('D41').Value = '2.782.93'
$ws.Range('E41').Value = '  +2.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.784'
$ws.Range('E42').Value = '  -2.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.44'
$ws.Range('E43').Value = '  -3.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.17'
$ws.Range('E44').Value = '  -4.34%  '
$ws.Range('E45').Value = '  -2.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.16'
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.08'
$ws.Range('E47').Value = '  -5.57%  '
$ws.Range('E48').Value = '  -7.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '319.34'
$ws.Range('E49').Value = '  -7.44%  '
$ws.Range('E50').Value = '  -3.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.979'
$ws.Range('E51').Value = '  -3.13%  '
